# Update "想去人数" (number of people interested) counts that changed
# between data refreshes for two rows that appear on both the "展览"
# sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览": row 4 (827 -> 828), row 5 (864 -> 866)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 828
$wsExhibit.Range("F5").Value = 866

# Sheet "全部类型": row 5 (827 -> 828), row 6 (864 -> 866)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 828
$wsAll.Range("F6").Value = 866
